# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets mirror the same data, so the same row/value updates are
# applied to each of them.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1093
    6  = 53
    8  = 11231
    9  = 4290
    13 = 2501
    14 = 1073
    15 = 107
    17 = 161
    19 = 11240
    20 = 11090
    25 = 35
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
